$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.240.39"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").Value = "1.830.67"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("D4").Value = "'1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'236.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").Value = "'0.6078"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.52%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "'0.07120"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.87%  "

$ws.Range("D9").Value = "'0.2818"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'23.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.84%  "

$ws.Range("D11").Value = "'0.07673"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").Value = "1.813.71"
$ws.Range("E12").Value = "  -1.81%  "

$ws.Range("D13").Value = "'4.834"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.91%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.00001011"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.6381"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.83%  "

$ws.Range("D16").Value = "2.080.77"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").Value = "'79.53"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'5.921"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.93%  "

$ws.Range("D19").Value = "29.225.38"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").Value = "'229.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("E21").Value = "  -3.85%  "

$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "'7.039"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.67%  "

$ws.Range("D25").Value = "'154.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.27%  "

$ws.Range("D26").Value = "'8.102"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.02%  "

$ws.Range("D27").Value = "'0.1295"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.51%  "

$ws.Range("D28").Value = "'16.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.61%  "

$ws.Range("D29").Value = "'1.490"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.43%  "

$ws.Range("D30").Value = "'0.06502"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.61%  "

$ws.Range("D31").Value = "'1.461"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.74%  "

$ws.Range("D32").Value = "'3.838"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.56%  "

$ws.Range("D33").Value = "'3.839"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.71%  "

$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").Value = "'1.746"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.69%  "

$ws.Range("D36").Value = "'0.6544"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.52%  "

$ws.Range("D37").Value = "'2.561"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.84%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.224.25"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.759"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.13%  "

$ws.Range("E40").Value = "  -4.96%  "

$ws.Range("D41").Value = "'6.536"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.62%  "

$ws.Range("D42").Value = "'0.9320"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.24%  "

$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").Value = "'101.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("D45").Value = "1.980.24"
$ws.Range("E45").Value = "  -1.18%  "

$ws.Range("D46").Value = "'63.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.03%  "

$ws.Range("D47").Value = "'0.00000000119"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D48").Value = "'1.613"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.03%  "

$ws.Range("D49").Value = "'8.558"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.66%  "

$ws.Range("D50").Value = "'6.502"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.75%  "

$ws.Range("D51").Value = "'0.05543"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.36%  "

